$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to reflect repulled data / recalculated mean
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = 4
$ws.Range("F5").Value = 2
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -11
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -1
$ws.Range("F14").Value = 6
$ws.Range("F15").Value = -3
